$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the first sheet "Export this as TSV" -> "Export as TSV"
$ws.Name = "Export as TSV"

# Freeze the header row (row 1) on this sheet.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

# Add error title / error message text to the existing data validations.
$dv = $ws.Range("I2:I1048576").Validation
$dv.ErrorTitle = "Value must come from list"
$dv.ErrorMessage = "Value must be one of: sequence."

$dv = $ws.Range("J2:J1048576").Validation
$dv.ErrorTitle = "Value must come from list"
$dv.ErrorMessage = "Value must be one of: scRNAseq-10xGenomics / scRNAseq / sciRNAseq / snRNAseq / SNARE2-RNAseq."

$dv = $ws.Range("K2:K1048576").Validation
$dv.ErrorTitle = "Value must come from list"
$dv.ErrorMessage = "Value must be one of: RNA."

$dv = $ws.Range("L2:L1048576").Validation
$dv.ErrorTitle = "Not a boolean"
$dv.ErrorMessage = 'The values in this column must be "TRUE" or "FALSE".'

$dv = $ws.Range("AA2:AA1048576").Validation
$dv.ErrorTitle = "Not a boolean"
$dv.ErrorMessage = 'The values in this column must be "TRUE" or "FALSE".'

$dv = $ws.Range("AG2:AG1048576").Validation
$dv.ErrorTitle = "Not a number"
$dv.ErrorMessage = "The values in this column must be numbers."

$dv = $ws.Range("AH2:AH1048576").Validation
$dv.ErrorTitle = "Value must come from list"
$dv.ErrorMessage = "Value must be one of: ng."

$dv = $ws.Range("AL2:AL1048576").Validation
$dv.ErrorTitle = "Not a number"
$dv.ErrorMessage = "The values in this column must be numbers."

$dv = $ws.Range("AM2:AM1048576").Validation
$dv.ErrorTitle = "Not a number"
$dv.ErrorMessage = "The values in this column must be numbers."
